# Update crypto price/volume table cells per the commit diff (data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.155.02"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "2.110.96"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.75"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.24"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "2.412.08"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.66"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.72"
$ws.Range("E14").Value = "  +4.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.786"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "2.114.32"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "38.055.07"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.04"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.87"
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.23"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.140"
$ws.Range("E27").Value = "  +11.35%  "
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.58"
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.70"
$ws.Range("E32").Value = "  +5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.61"
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.48"
$ws.Range("E36").Value = "  +6.69%  "
$ws.Range("E37").Value = "  +4.99%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("E39").Value = "  -3.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0994"
$ws.Range("E40").Value = "  +7.01%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.48"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0216"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").Value = "1.462.90"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.13"
$ws.Range("E46").Value = "  +7.31%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.17"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.03"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.33"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").Value = "2.307.83"
$ws.Range("E51").Value = "  +2.85%  "
